$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BA")

# --- Re-format the table: the whole range was stored as Text ("@"); switch the
# header + text columns to General and the numeric columns to a plain Number
# format, which is what lets the new values below store as real numbers. ---

# Header row - keep its bold styling, switch off the Text format.
$ws.Range("A1:K1").NumberFormat = "general"

# Text columns in the body (BA_Name, state_name, is_msmed, contact_person_name).
$ws.Range("B2:C4").NumberFormat = "general"
$ws.Range("H2:H4").NumberFormat = "general"
$ws.Range("J2:J4").NumberFormat = "general"

# Hyperlink column (email_id) - keep the hyperlink styling, General format.
$ws.Range("I2:I4").NumberFormat = "general"

# Numeric columns get a real number format.
$ws.Range("A2:A4").NumberFormat = "0"
$ws.Range("D2:G4").NumberFormat = "0"
$ws.Range("K2:K4").NumberFormat = "0"

# --- Update the data values ---
$ws.Range("A2").Value2 = 432
$ws.Range("K2").Value2 = 9876543210
$ws.Range("B3").Value2 = "Nish"
$ws.Range("A4").Value2 = 431

# --- Add a new (mostly blank) row 5, formatted like the other numeric cells ---
$ws.Range("A5").NumberFormat = "0"
$ws.Range("K5").NumberFormat = "0"

# --- Update the active selection ---
$ws.Range("B3").Select()
